$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 74; existing rows 74:155 shift down to 75:156.
$ws.Rows("74:74").Insert()

# Populate the newly inserted row 74 with the new weekly record
# (context columns match the surrounding "Granada" / Vega Modelo de Temuco block).
$ws.Range("A74").Value = 10
$ws.Range("B74").Value = "Vega Modelo de Temuco"
$ws.Range("C74").Value = "La Araucanía"
$ws.Range("D74").Value = 44789
$ws.Range("E74").Value = 9
$ws.Range("F74").Value = "Fruta"
$ws.Range("G74").Value = 100104
$ws.Range("H74").Value = "Frutos de pepita"
$ws.Range("I74").Value = 100104001
$ws.Range("J74").Value = "Granada"
$ws.Range("K74").Value = "Wonderfull"
$ws.Range("L74").Value = "Primera"
$ws.Range("M74").Value = 80
$ws.Range("N74").Value = 14000
$ws.Range("O74").Value = 14000
$ws.Range("P74").Value = 14000
$ws.Range("Q74").Value = "$/bandeja 10 kilos granel"
$ws.Range("R74").Value = "Provincia de Limarí"
$ws.Range("S74").Value = 1400
$ws.Range("T74").Value = 10
